$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 19.239999999999998
$ws.Range("E2").Value = 18.670000000000002
$ws.Range("F2").Value = -2.9625779625779458
$ws.Range("G2").Value = $false

# Row 3 new data
$ws.Range("C3").Value = 9703.74
